$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(21, 8).Value = 26803.8
$ws.Cells.Item(21, 10).Value = 21000
$ws.Cells.Item(21, 12).Value = 21000
$ws.Cells.Item(21, 14).Value = -21936
$ws.Cells.Item(23, 8).Value = 26803.8
$ws.Cells.Item(23, 10).Value = 21000
$ws.Cells.Item(23, 12).Value = 21000
$ws.Cells.Item(23, 14).Value = -21468
$ws.Cells.Item(64, 8).Value = 3242.6316
$ws.Cells.Item(64, 9).Value = 3603.75
$ws.Cells.Item(64, 10).Value = 2980
$ws.Cells.Item(64, 11).Value = 3603.75
$ws.Cells.Item(64, 12).Value = 2980
$ws.Cells.Item(64, 13).Value = -3355.75
$ws.Cells.Item(64, 14).Value = -3476
$ws.Cells.Item(67, 8).Value = 3242.6316
$ws.Cells.Item(67, 9).Value = 3603.75
$ws.Cells.Item(67, 10).Value = 2980
$ws.Cells.Item(67, 11).Value = 3603.75
$ws.Cells.Item(67, 12).Value = 2980
$ws.Cells.Item(67, 13).Value = -2745.75
$ws.Cells.Item(67, 14).Value = -4696
$ws.Cells.Item(69, 8).Value = 4550
$ws.Cells.Item(69, 10).Value = 4550
$ws.Cells.Item(69, 12).Value = 13650
$ws.Cells.Item(69, 14).Value = -15398
$ws.Cells.Item(72, 8).Value = 4550
$ws.Cells.Item(72, 10).Value = 4550
$ws.Cells.Item(72, 12).Value = 40950
$ws.Cells.Item(72, 14).Value = -49686
$ws.Cells.Item(76, 8).Value = 3802.5
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 10).Value = 3802.5
$ws.Cells.Item(76, 11).Value = 0
$ws.Cells.Item(76, 12).Value = 3802.5
$ws.Cells.Item(76, 13).ClearContents()
$ws.Cells.Item(76, 14).Value = -4432.5
$ws.Cells.Item(79, 8).Value = 3802.5
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 10).Value = 3802.5
$ws.Cells.Item(79, 11).Value = 0
$ws.Cells.Item(79, 12).Value = 3802.5
$ws.Cells.Item(79, 13).ClearContents()
$ws.Cells.Item(79, 14).Value = -5986.5
$ws.Cells.Item(132, 8).Value = 2470959.2
$ws.Cells.Item(132, 9).Value = 2668024
$ws.Cells.Item(132, 11).Value = 8004072
$ws.Cells.Item(132, 13).Value = -8001542
$ws.Cells.Item(135, 8).Value = 983.0741
$ws.Cells.Item(135, 9).Value = 524.6111
$ws.Cells.Item(135, 11).Value = 4721.4999
$ws.Cells.Item(135, 13).Value = -2186.4999
$ws.Cells.Item(138, 8).Value = 5985.473
$ws.Cells.Item(138, 9).Value = 2688.4482
$ws.Cells.Item(138, 10).Value = 8110.222
$ws.Cells.Item(138, 11).Value = 8065.344599999999
$ws.Cells.Item(138, 12).Value = 24330.666
$ws.Cells.Item(138, 13).Value = -2925.344599999999
$ws.Cells.Item(138, 14).Value = -34610.666
$ws.Cells.Item(141, 8).Value = 1545522.2
$ws.Cells.Item(141, 9).Value = 1999.75
$ws.Cells.Item(141, 10).Value = 3089044.8
$ws.Cells.Item(141, 11).Value = 5999.25
$ws.Cells.Item(141, 12).Value = 9267134.399999999
$ws.Cells.Item(141, 13).Value = -819.25
$ws.Cells.Item(141, 14).Value = -9277494.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 21381.072
$ws.Cells.Item(32, 9).Value = 17346.242
$ws.Cells.Item(32, 11).Value = 17346.242
$ws.Cells.Item(32, 13).Value = -17059.242
$ws.Cells.Item(80, 8).Value = 29649.143
$ws.Cells.Item(80, 10).Value = 29649.143
$ws.Cells.Item(80, 12).Value = 29649.143
$ws.Cells.Item(80, 14).Value = -31645.143
$ws.Cells.Item(83, 8).Value = 29649.143
$ws.Cells.Item(83, 10).Value = 29649.143
$ws.Cells.Item(83, 12).Value = 88947.429
$ws.Cells.Item(83, 14).Value = -98931.429
$ws.Cells.Item(110, 8).Value = 3146.625
$ws.Cells.Item(110, 9).Value = 636.41174
$ws.Cells.Item(110, 10).Value = 9242.857
$ws.Cells.Item(110, 11).Value = 636.41174
$ws.Cells.Item(110, 12).Value = 9242.857
$ws.Cells.Item(110, 13).Value = 1408.58826
$ws.Cells.Item(110, 14).Value = -13332.857

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(82, 8).Value = 23017.334
$ws.Cells.Item(82, 10).Value = 29045.445
$ws.Cells.Item(82, 12).Value = 29045.445
$ws.Cells.Item(82, 14).Value = -29811.445
$ws.Cells.Item(85, 8).Value = 23017.334
$ws.Cells.Item(85, 10).Value = 29045.445
$ws.Cells.Item(85, 12).Value = 29045.445
$ws.Cells.Item(85, 14).Value = -31697.445
$ws.Cells.Item(94, 8).Value = 751.19354
$ws.Cells.Item(94, 9).Value = 588.0454999999999
$ws.Cells.Item(94, 10).Value = 1150
$ws.Cells.Item(94, 11).Value = 588.0454999999999
$ws.Cells.Item(94, 12).Value = 1150
$ws.Cells.Item(94, 13).Value = -137.0454999999999
$ws.Cells.Item(94, 14).Value = -2052
$ws.Cells.Item(105, 8).Value = 2395.8948
$ws.Cells.Item(105, 9).Value = 2146.8462
$ws.Cells.Item(105, 10).Value = 2935.5
$ws.Cells.Item(105, 11).Value = 2146.8462
$ws.Cells.Item(105, 12).Value = 2935.5
$ws.Cells.Item(105, 13).Value = -399.8462
$ws.Cells.Item(105, 14).Value = -6429.5
$ws.Cells.Item(114, 8).Value = 0
$ws.Cells.Item(114, 10).Value = 0
$ws.Cells.Item(114, 12).Value = 0
$ws.Cells.Item(114, 14).ClearContents()
$ws.Cells.Item(134, 8).Value = 2984.04
$ws.Cells.Item(134, 9).Value = 2679.8
$ws.Cells.Item(134, 11).Value = 8039.400000000001
$ws.Cells.Item(134, 13).Value = -5504.400000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 3166.6667
$ws.Cells.Item(16, 9).Value = 1000
$ws.Cells.Item(16, 10).Value = 3600
$ws.Cells.Item(16, 11).Value = 1000
$ws.Cells.Item(16, 12).Value = 3600
$ws.Cells.Item(16, 13).Value = -713
$ws.Cells.Item(16, 14).Value = -4174
$ws.Cells.Item(31, 8).Value = 1518107.9
$ws.Cells.Item(31, 9).Value = 2175594.5
$ws.Cells.Item(31, 11).Value = 2175594.5
$ws.Cells.Item(31, 13).Value = -2175299.5
$ws.Cells.Item(34, 8).Value = 1518107.9
$ws.Cells.Item(34, 9).Value = 2175594.5
$ws.Cells.Item(34, 11).Value = 2175594.5
$ws.Cells.Item(34, 13).Value = -2175392.5
$ws.Cells.Item(58, 8).Value = 13517619
$ws.Cells.Item(58, 9).Value = 2290.2222
$ws.Cells.Item(58, 10).Value = 26321616
$ws.Cells.Item(58, 11).Value = 2290.2222
$ws.Cells.Item(58, 12).Value = 26321616
$ws.Cells.Item(58, 13).Value = -2087.2222
$ws.Cells.Item(58, 14).Value = -26322022
$ws.Cells.Item(74, 8).Value = 24332.666
$ws.Cells.Item(74, 10).Value = 24332.666
$ws.Cells.Item(74, 12).Value = 24332.666
$ws.Cells.Item(74, 14).Value = -26080.666
$ws.Cells.Item(77, 8).Value = 24332.666
$ws.Cells.Item(77, 10).Value = 24332.666
$ws.Cells.Item(77, 12).Value = 72997.99800000001
$ws.Cells.Item(77, 14).Value = -81733.99800000001
$ws.Cells.Item(113, 8).Value = 3166.6667
$ws.Cells.Item(113, 9).Value = 1000
$ws.Cells.Item(113, 10).Value = 3600
$ws.Cells.Item(113, 11).Value = 1000
$ws.Cells.Item(113, 12).Value = 3600
$ws.Cells.Item(113, 13).Value = 1170
$ws.Cells.Item(113, 14).Value = -7940
$ws.Cells.Item(136, 8).Value = 13517619
$ws.Cells.Item(136, 9).Value = 2290.2222
$ws.Cells.Item(136, 10).Value = 26321616
$ws.Cells.Item(136, 11).Value = 6870.6666
$ws.Cells.Item(136, 12).Value = 78964848
$ws.Cells.Item(136, 13).Value = -4320.6666
$ws.Cells.Item(136, 14).Value = -78969948

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 51.666668
$ws.Cells.Item(2, 9).Value = 25
$ws.Cells.Item(2, 10).Value = 69.44444
$ws.Cells.Item(2, 11).Value = 150
$ws.Cells.Item(2, 12).Value = 416.66664
$ws.Cells.Item(2, 13).Value = -37
$ws.Cells.Item(2, 14).Value = -642.66664
$ws.Cells.Item(5, 8).Value = 792.05554
$ws.Cells.Item(5, 10).Value = 3100
$ws.Cells.Item(5, 12).Value = 9300
$ws.Cells.Item(5, 14).Value = -9524
$ws.Cells.Item(92, 8).Value = 1841.9166
$ws.Cells.Item(92, 9).Value = 900
$ws.Cells.Item(92, 10).Value = 1927.5454
$ws.Cells.Item(92, 11).Value = 2700
$ws.Cells.Item(92, 12).Value = 5782.6362
$ws.Cells.Item(92, 13).Value = -1452
$ws.Cells.Item(92, 14).Value = -8278.636200000001
$ws.Cells.Item(135, 8).Value = 792.05554
$ws.Cells.Item(135, 10).Value = 3100
$ws.Cells.Item(135, 12).Value = 27900
$ws.Cells.Item(135, 14).Value = -32970

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4293.852
$ws.Cells.Item(70, 9).Value = 4406.2666
$ws.Cells.Item(70, 10).Value = 4153.3335
$ws.Cells.Item(70, 11).Value = 4406.2666
$ws.Cells.Item(70, 12).Value = 4153.3335
$ws.Cells.Item(70, 13).Value = -4136.2666
$ws.Cells.Item(70, 14).Value = -4693.3335
$ws.Cells.Item(73, 8).Value = 4293.852
$ws.Cells.Item(73, 9).Value = 4406.2666
$ws.Cells.Item(73, 10).Value = 4153.3335
$ws.Cells.Item(73, 11).Value = 4406.2666
$ws.Cells.Item(73, 12).Value = 4153.3335
$ws.Cells.Item(73, 13).Value = -3470.2666
$ws.Cells.Item(73, 14).Value = -6025.3335
$ws.Cells.Item(80, 8).Value = 3512.8572
$ws.Cells.Item(80, 9).Value = 3573.3333
$ws.Cells.Item(80, 10).Value = 3150
$ws.Cells.Item(80, 11).Value = 3573.3333
$ws.Cells.Item(80, 12).Value = 3150
$ws.Cells.Item(80, 13).Value = -2575.3333
$ws.Cells.Item(80, 14).Value = -5146
$ws.Cells.Item(83, 8).Value = 3512.8572
$ws.Cells.Item(83, 9).Value = 3573.3333
$ws.Cells.Item(83, 10).Value = 3150
$ws.Cells.Item(83, 11).Value = 17866.6665
$ws.Cells.Item(83, 12).Value = 15750
$ws.Cells.Item(83, 13).Value = -12874.6665
$ws.Cells.Item(83, 14).Value = -25734

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 2862577.2
$ws.Cells.Item(136, 9).Value = 4766447.5
$ws.Cells.Item(136, 10).Value = 6772.143
$ws.Cells.Item(136, 11).Value = 14299342.5
$ws.Cells.Item(136, 12).Value = 20316.429
$ws.Cells.Item(136, 13).Value = -14296792.5
$ws.Cells.Item(136, 14).Value = -25416.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2646690.8
$ws.Cells.Item(132, 9).Value = 3230324
$ws.Cells.Item(132, 11).Value = 9690972
$ws.Cells.Item(132, 13).Value = -9688442
$ws.Cells.Item(136, 8).Value = 2134.83
$ws.Cells.Item(136, 9).Value = 1586.125
$ws.Cells.Item(136, 10).Value = 3823.1538
$ws.Cells.Item(136, 11).Value = 4758.375
$ws.Cells.Item(136, 12).Value = 11469.4614
$ws.Cells.Item(136, 13).Value = -2208.375
$ws.Cells.Item(136, 14).Value = -16569.4614
